$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G1:H1").EntireColumn.Delete()
